$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Auth0" (the trigger program, cell B14) is renamed to describe the
# heartbeat generation/reception program instead.
$ws.Range("B14").Value = "Generación y recepción de HeartBeats"

# Row 5 ("Nicolas Simmonds") no longer contributes to the "Programa
# Disparador" effort column, row 6 ("Luis Carlos Garavito") now does.
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 33.333

# Column B is widened to fit the longer artifact description that now
# lives in it.
$ws.Columns("B").ColumnWidth = 53

# Reflect where the user's cursor ended up after the edits.
$ws.Range("B15").Select()

# Page is set up for portrait printing.
$ws.PageSetup.Orientation = 1
